{"js": "const replacements = [\n  ['2024-06-23 Sunday', '2024-06-24 Monday'],\n  ['748\u00f73=249, 1', '616\u00f73=205, 1'],\n  ['486\u00f79=54, 0', '832\u00f78=104, 0'],\n  ['215\u00f73=71, 2', '165\u00f72=82, 1'],\n  ['766\u00f76=127, 4', '406\u00f73=135, 1'],\n  ['600\u00f76=100, 0', '852\u00f78=106, 4'],\n  ['573\u00f79=63, 6', '862\u00f73=287, 1'],\n  ['309\u00f74=77, 1', '137\u00f78=17, 1'],\n  ['344\u00f79=38, 2', '498\u00f73=166, 0'],\n  ['649\u00f74=162, 1', '853\u00f79=94, 7'],\n  ['501\u00f78=62, 5', '206\u00f76=34, 2'],\n  ['397\u00f72=198, 1', '830\u00f75=166, 0'],\n  ['487\u00f78=60, 7', '606\u00f74=151, 2'],\n  ['411\u00f78=51, 3', '841\u00f73=280, 1'],\n  ['415\u00f75=83, 0', '725\u00f77=103, 4'],\n  ['195\u00f74=48, 3', '993\u00f75=198, 3'],\n  ['881\u00f73=293, 2', '631\u00f77=90, 1'],\n  ['786\u00f75=157, 1', '976\u00f76=162, 4'],\n  ['510\u00f76=85, 0', '411\u00f79=45, 6'],\n  ['975\u00f75=195, 0', '485\u00f74=121, 1'],\n  ['285\u00f74=71, 1', '331\u00f72=165, 1'],\n  ['137\u00f72=68, 1', '627\u00f75=125, 2'],\n  ['121\u00f72=60, 1', '572\u00f73=190, 2'],\n  ['586\u00f72=293, 0', '107\u00f72=53, 1'],\n  ['557\u00f72=278, 1', '901\u00f73=300, 1'],\n  ['283\u00f75=56, 3', '507\u00f73=169, 0'],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-06-23 Sunday', '2024-06-24 Monday'),\n    @('748\u00f73=249, 1', '616\u00f73=205, 1'),\n    @('486\u00f79=54, 0', '832\u00f78=104, 0'),\n    @('215\u00f73=71, 2', '165\u00f72=82, 1'),\n    @('766\u00f76=127, 4', '406\u00f73=135, 1'),\n    @('600\u00f76=100, 0', '852\u00f78=106, 4'),\n    @('573\u00f79=63, 6', '862\u00f73=287, 1'),\n    @('309\u00f74=77, 1', '137\u00f78=17, 1'),\n    @('344\u00f79=38, 2', '498\u00f73=166, 0'),\n    @('649\u00f74=162, 1', '853\u00f79=94, 7'),\n    @('501\u00f78=62, 5', '206\u00f76=34, 2'),\n    @('397\u00f72=198, 1', '830\u00f75=166, 0'),\n    @('487\u00f78=60, 7', '606\u00f74=151, 2'),\n    @('411\u00f78=51, 3', '841\u00f73=280, 1'),\n    @('415\u00f75=83, 0', '725\u00f77=103, 4'),\n    @('195\u00f74=48, 3', '993\u00f75=198, 3'),\n    @('881\u00f73=293, 2', '631\u00f77=90, 1'),\n    @('786\u00f75=157, 1', '976\u00f76=162, 4'),\n    @('510\u00f76=85, 0', '411\u00f79=45, 6'),\n    @('975\u00f75=195, 0', '485\u00f74=121, 1'),\n    @('285\u00f74=71, 1', '331\u00f72=165, 1'),\n    @('137\u00f72=68, 1', '627\u00f75=125, 2'),\n    @('121\u00f72=60, 1', '572\u00f73=190, 2'),\n    @('586\u00f72=293, 0', '107\u00f72=53, 1'),\n    @('557\u00f72=278, 1', '901\u00f73=300, 1'),\n    @('283\u00f75=56, 3', '507\u00f73=169, 0'),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null  # wdReplaceAll\n}"}
